$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.893.96"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "1.806.99"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4649"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3714"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07378"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8743"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("D12").Value = "1.870.40"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.374"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.84"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.495"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07034"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008721"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.68"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").Value = "26.898.41"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.307"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "2.007.55"
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.904"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.50"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.149"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.288"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08927"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7598"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.158"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.460"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.913"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.106"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05257"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.928"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.242"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.378"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5293"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1665"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.529"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4993"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.32"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.95"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9999"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06297"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.57%  "
